$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1: copy H1's format (bold, bordered, centered style)
# then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-27.
$data = @(
    @(4, 4),
    @(5, 7),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 3),
    @(1, 6),
    @(1, 7),
    @(4, 7),
    @(1, 7),
    @(1, 7),
    @(1, 8),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 7),
    @(1, 5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
